$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.89"
$ws.Range("D3").Value = "'26.08"
$ws.Range("D4").Value = "'5.081"
$ws.Range("D5").Value = "'0.05594"
$ws.Range("D6").Value = "'6.478"
$ws.Range("D7").Value = "'3.023"
$ws.Range("D8").Value = "'0.8120"
$ws.Range("D9").Value = "'0.8431"
$ws.Range("D10").Value = "'0.1344"
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = "'0.02770"
$ws.Range("E11").Value = '10BitrueCoinBTR'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = "'0.09400"
$ws.Range("E12").Value = '11BitMartTokenBMX'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = "'0.001516"
$ws.Range("E13").Value = '12BitForexTokenBF'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = "'0.0006030"
$ws.Range("E14").Value = '13OneONE'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = "'0.006158"
$ws.Range("E15").Value = '14TigerCashTCH'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = "'3.557"
$ws.Range("E16").Value = '15LEOLEO'
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D17").Value = "'2.118"
$ws.Range("E17").Value = '16BTSETokenBTSE'
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D18").Value = "'0.3183"
$ws.Range("E18").Value = '17BitpandaEcosystemTokenBEST'
$ws.Range("B19").Value = 'MandalaExchangeToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D19").Value = "'0.07006"
$ws.Range("E19").Value = '18MandalaExchangeTokenMDX'
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = "'0.03225"
$ws.Range("E20").Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("D21").Value = "'0.1320"
$ws.Range("D22").Value = "'3.762"
$ws.Range("D23").Value = "'0.04700"
$ws.Range("D24").Value = "'0.1375"
$ws.Range("D25").Value = "'0.001249"
$ws.Range("D26").Value = "'0.004613"
$ws.Range("D27").Value = "'0.00009600"
$ws.Range("D28").Value = "'0.0001390"
$ws.Range("D40").Value = "'0.03659"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.006120"
$ws.Range("E41").Value = '40KickTokenKICKBestin24h'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1053"
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("D43").Value = "'0.002500"
$ws.Range("E43").Value = '42CEJICEJIWorstin24h'
$ws.Range("D44").Value = "'0.008712"
$ws.Range("D45").Value = "'0.00005296"
$ws.Range("D47").Value = "'0.1328"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
